$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.534.21"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D5").Value = "'228.65"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'34.98"
$ws.Range("E8").Value = "  +6.36%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "2.074.40"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "'11.23"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "1.818.91"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "'0.650"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("D17").Value = "34.513.48"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'69.32"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'246.44"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "'11.47"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'172.52"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").Value = "'8.12"
$ws.Range("E26").Value = "  +10.66%  "
$ws.Range("D27").Value = "'16.84"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'4.04"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'0.0537"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "1.398.39"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").Value = "'0.680"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "'83.87"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").Value = "'0.966"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.14"
$ws.Range("E44").Value = "  +5.70%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'13.30"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "'6.01"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "1.973.65"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "'105.48"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("E51").Value = "  +0.10%  "